$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "天际股份"

$ws.Range("A3").Value = "天际股份"
$ws.Range("B3").Value = "天际股份"

$ws.Range("A4").Value = "平潭发展"
$ws.Range("C4").Value = "世纪华通"

$ws.Range("A5").Value = "方正电机"
$ws.Range("B5").Value = "农业银行"
$ws.Range("C5").Value = "农业银行"

$ws.Range("A6").Value = "东百集团"
$ws.Range("B6").Value = "福龙马"
$ws.Range("C6").Value = "平潭发展"

$ws.Range("A7").Value = "孚日股份"
$ws.Range("B7").Value = "摩恩电气"
$ws.Range("C7").Value = "安泰集团"

$ws.Range("A8").Value = "农业银行"
$ws.Range("B8").Value = "协鑫集成"
$ws.Range("C8").Value = "摩恩电气"

$ws.Range("A9").Value = "隆基绿能"
$ws.Range("B9").Value = "方正电机"
$ws.Range("C9").Value = "孚日股份"

$ws.Range("A10").Value = "摩恩电气"
$ws.Range("B10").Value = "东百集团"
$ws.Range("C10").Value = "海马汽车"

$ws.Range("A11").Value = "香农芯创"
$ws.Range("B11").Value = "香农芯创"
$ws.Range("C11").Value = "香农芯创"

$ws.Range("A12").Value = "安泰集团"
$ws.Range("C12").Value = "多氟多"

$ws.Range("A13").Value = "福龙马"
$ws.Range("B13").Value = "安泰集团"
$ws.Range("C13").Value = "方正电机"

$ws.Range("A14").Value = "多氟多"
$ws.Range("B14").Value = "中国中免"
$ws.Range("C14").Value = "东百集团"

$ws.Range("A15").Value = "海马汽车"
$ws.Range("B15").Value = "隆基绿能"
$ws.Range("C15").Value = "隆基绿能"

$ws.Range("A16").Value = "人民同泰"
$ws.Range("B16").Value = "永太科技"
$ws.Range("C16").Value = "上海电力"

$ws.Range("A17").Value = "协鑫集成"
$ws.Range("B17").Value = "中利集团"
$ws.Range("C17").Value = "万向钱潮"

$ws.Range("A18").Value = "永太科技"
$ws.Range("B18").Value = "特变电工"
$ws.Range("C18").Value = "航天智装"

$ws.Range("A19").Value = "石大胜华"
$ws.Range("B19").Value = "孚日股份"
$ws.Range("C19").Value = "福龙马"

$ws.Range("A20").Value = "世纪华通"
$ws.Range("B20").Value = "众生药业"
$ws.Range("C20").Value = "天下秀"

$ws.Range("A21").Value = "众生药业"
$ws.Range("B21").Value = "人民同泰"
$ws.Range("C21").Value = "上海沪工"
